# Scheduled runner: refresh cached Universalis price snapshots + dependent
# profit calcs (currentAveragePrice*, LevePrice*, LeveProfit*) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2985.3538
$ws.Range("I15").Value = 2985.3538
$ws.Range("K15").Value = 8956.061399999999
$ws.Range("M15").Value = -8787.061399999999
$ws.Range("H51").Value = 1838.6666
$ws.Range("I51").Value = 863.6667
$ws.Range("J51").Value = 2082.4167
$ws.Range("K51").Value = 863.6667
$ws.Range("L51").Value = 2082.4167
$ws.Range("M51").Value = -379.6667
$ws.Range("N51").Value = -3050.4167
$ws.Range("H62").Value = 13894076
$ws.Range("I62").Value = 18523102
$ws.Range("J62").Value = 6999.5
$ws.Range("K62").Value = 18523102
$ws.Range("L62").Value = 6999.5
$ws.Range("M62").Value = -18522478
$ws.Range("N62").Value = -8247.5
$ws.Range("H65").Value = 13894076
$ws.Range("I65").Value = 18523102
$ws.Range("J65").Value = 6999.5
$ws.Range("K65").Value = 92615510
$ws.Range("L65").Value = 34997.5
$ws.Range("M65").Value = -92612390
$ws.Range("N65").Value = -41237.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 97
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 97
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 97
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -329
$ws.Range("H61").Value = 20834226
$ws.Range("I61").Value = 22223022
$ws.Range("J61").Value = 2271.3333
$ws.Range("K61").Value = 22223022
$ws.Range("L61").Value = 2271.3333
$ws.Range("M61").Value = -22222810
$ws.Range("N61").Value = -2695.3333
$ws.Range("H63").Value = 66668800
$ws.Range("I63").Value = 1966.125
$ws.Range("K63").Value = 1966.125
$ws.Range("M63").Value = -1280.125
$ws.Range("H66").Value = 66668800
$ws.Range("I66").Value = 1966.125
$ws.Range("K66").Value = 9830.625
$ws.Range("M66").Value = -6398.625
$ws.Range("H74").Value = 963.9211
$ws.Range("I74").Value = 723.9706
$ws.Range("J74").Value = 3003.5
$ws.Range("K74").Value = 723.9706
$ws.Range("L74").Value = 3003.5
$ws.Range("M74").Value = 150.0294
$ws.Range("N74").Value = -4751.5
$ws.Range("H77").Value = 963.9211
$ws.Range("I77").Value = 723.9706
$ws.Range("J77").Value = 3003.5
$ws.Range("K77").Value = 3619.853
$ws.Range("L77").Value = 15017.5
$ws.Range("M77").Value = 748.1469999999999
$ws.Range("N77").Value = -23753.5
$ws.Range("H110").Value = 1473.28
$ws.Range("I110").Value = 1065.7646
$ws.Range("K110").Value = 1065.7646
$ws.Range("M110").Value = 979.2354
$ws.Range("H122").Value = 3214.2856
$ws.Range("I122").Value = 3166.6667
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9500.000100000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -7050.000100000001
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 2742.5293
$ws.Range("I132").Value = 2743.0833
$ws.Range("J132").Value = 2741.2
$ws.Range("K132").Value = 8229.249899999999
$ws.Range("L132").Value = 8223.599999999999
$ws.Range("M132").Value = -5699.249899999999
$ws.Range("N132").Value = -13283.6
$ws.Range("H136").Value = 20834226
$ws.Range("I136").Value = 22223022
$ws.Range("J136").Value = 2271.3333
$ws.Range("K136").Value = 66669066
$ws.Range("L136").Value = 6813.999899999999
$ws.Range("M136").Value = -66666516
$ws.Range("N136").Value = -11913.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2664.7856
$ws.Range("I20").Value = 1977.4445
$ws.Range("K20").Value = 1977.4445
$ws.Range("M20").Value = -1730.4445
$ws.Range("H86").Value = 3266.5715
$ws.Range("I86").Value = 3722.85
$ws.Range("J86").Value = 2125.875
$ws.Range("K86").Value = 3722.85
$ws.Range("L86").Value = 2125.875
$ws.Range("M86").Value = -2599.85
$ws.Range("N86").Value = -4371.875
$ws.Range("H89").Value = 3266.5715
$ws.Range("I89").Value = 3722.85
$ws.Range("J89").Value = 2125.875
$ws.Range("K89").Value = 18614.25
$ws.Range("L89").Value = 10629.375
$ws.Range("M89").Value = -12998.25
$ws.Range("N89").Value = -21861.375
$ws.Range("H107").Value = 2191
$ws.Range("I107").Value = 1633.1666
$ws.Range("J107").Value = 3027.75
$ws.Range("K107").Value = 1633.1666
$ws.Range("L107").Value = 3027.75
$ws.Range("M107").Value = 286.8334
$ws.Range("N107").Value = -6867.75
$ws.Range("H134").Value = 5149.8716
$ws.Range("I134").Value = 1959.6897
$ws.Range("K134").Value = 5879.0691
$ws.Range("M134").Value = -3344.0691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1341.0667
$ws.Range("I122").Value = 1225.375
$ws.Range("K122").Value = 3676.125
$ws.Range("M122").Value = -1226.125
$ws.Range("H132").Value = 3125.7693
$ws.Range("I132").Value = 2954.7058
$ws.Range("J132").Value = 3448.889
$ws.Range("K132").Value = 8864.117400000001
$ws.Range("L132").Value = 10346.667
$ws.Range("M132").Value = -6334.117400000001
$ws.Range("N132").Value = -15406.667
$ws.Range("H134").Value = 26316978
$ws.Range("I134").Value = 1212.8667
$ws.Range("J134").Value = 125001096
$ws.Range("K134").Value = 3638.6001
$ws.Range("L134").Value = 375003288
$ws.Range("M134").Value = -1103.6001
$ws.Range("N134").Value = -375008358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 97.82143000000001
$ws.Range("I12").Value = 157.28572
$ws.Range("J12").Value = 78
$ws.Range("K12").Value = 471.85716
$ws.Range("L12").Value = 234
$ws.Range("M12").Value = -298.85716
$ws.Range("N12").Value = -580
$ws.Range("H13").Value = 539.6
$ws.Range("I13").Value = 233
$ws.Range("K13").Value = 699
$ws.Range("M13").Value = -531
$ws.Range("H131").Value = 18519758
$ws.Range("J131").Value = 1325.0204
$ws.Range("L131").Value = 3975.0612
$ws.Range("N131").Value = -14055.0612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2959.0908
$ws.Range("I80").Value = 1710
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 1710
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -712
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 2959.0908
$ws.Range("I83").Value = 1710
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 8550
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -3558
$ws.Range("N83").Value = -29984
$ws.Range("H132").Value = 1963.6072
$ws.Range("I132").Value = 1545.7894
$ws.Range("K132").Value = 4637.3682
$ws.Range("M132").Value = -2107.3682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7250
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 11500
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 11500
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -11772
$ws.Range("H122").Value = 17861488
$ws.Range("I122").Value = 25001342
$ws.Range("K122").Value = 75004026
$ws.Range("M122").Value = -75001576
$ws.Range("H132").Value = 19951.51
$ws.Range("I132").Value = 1225.625
$ws.Range("J132").Value = 46004.914
$ws.Range("K132").Value = 3676.875
$ws.Range("L132").Value = 138014.742
$ws.Range("M132").Value = -1146.875
$ws.Range("N132").Value = -143074.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2424.682
$ws.Range("I132").Value = 2508
$ws.Range("J132").Value = 2049.75
$ws.Range("K132").Value = 7524
$ws.Range("L132").Value = 6149.25
$ws.Range("M132").Value = -4994
$ws.Range("N132").Value = -11209.25
